$d = $word.ActiveDocument

# 1. Rework the "Programming Languages" line in the Skills section:
#    "Programming Languages: Proficient in Java and SQL, knowledgeable in
#     Python, C++, Racket, HTML + CSS"
#    becomes
#    "Languages: Proficient in Java and SQL, knowledgeable in Python,
#     HTML + CSS, learning React.js, Bootstrap"
$d.Content.Find.Execute("Programming Languages:", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Languages:", 2) | Out-Null

$d.Content.Find.Execute("knowledgeable in Python, C++, Racket, HTML + CSS", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "knowledgeable in Python, HTML + CSS, learning React.js, Bootstrap", 2) | Out-Null

# Move the "_GoBack" bookmark (currently sitting in front of the
# "Education" heading paragraph) into its new home: a zero-length range
# sitting right before "SQL" on that same Skills line. Re-adding a
# bookmark under a name that already exists relocates it, so this both
# removes the old one and places the new one in a single step.
$r = $d.Content
$r.Find.Execute("Java and SQL", $true, $false, $false, $false, $false, `
                 $true, 1, $false, "", 0) | Out-Null
$sqlStart = $r.End - 3
$bmRange = $d.Range($sqlStart, $sqlStart)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 2. Append ", Linux" to the end of the "Software and IDEs" line, as a new
#    run placed right after the existing "Git" run (which stays intact,
#    spell-check markers and all). Scope the search to that one paragraph
#    so we don't touch the later "Github" mention in the Work Experience
#    section.
$ideParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Software and IDEs*") {
        $ideParagraph = $p.Range
        break
    }
}
$ideParagraph.Find.Execute("Git", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0) | Out-Null
$gitEnd = $ideParagraph.End
$ideParagraph.InsertAfter(", Linux")
$linuxRange = $d.Range($gitEnd, $gitEnd + 7)
$linuxRange.Font.Name = "Segoe UI"
$linuxRange.Font.Size = 11
